$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 340, pushing the old rows
# 340-344 down to become rows 342-346 (formatting, e.g. the date style on
# column D, carries down automatically with the insert).
$ws.Rows("340:341").Insert()

# New row 340 data
$ws.Cells.Item(340, 1).Value = 7
$ws.Cells.Item(340, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(340, 3).Value = "Ñuble"
$ws.Cells.Item(340, 4).Value = 44890
$ws.Cells.Item(340, 5).Value = 16
$ws.Cells.Item(340, 6).Value = 100112023
$ws.Cells.Item(340, 7).Value = "Brócoli"
$ws.Cells.Item(340, 8).Value = "Sin especificar"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 240
$ws.Cells.Item(340, 11).Value = 550
$ws.Cells.Item(340, 12).Value = 600
$ws.Cells.Item(340, 13).Value = 575
$ws.Cells.Item(340, 14).Value = "$/unidad"
$ws.Cells.Item(340, 15).Value = "Región del Maule"
$ws.Cells.Item(340, 16).Value = 575
$ws.Cells.Item(340, 17).Value = 1
$ws.Cells.Item(340, 18).Value = "Hortaliza"

# New row 341 data
$ws.Cells.Item(341, 1).Value = 7
$ws.Cells.Item(341, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(341, 3).Value = "Ñuble"
$ws.Cells.Item(341, 4).Value = 44890
$ws.Cells.Item(341, 5).Value = 16
$ws.Cells.Item(341, 6).Value = 100112023
$ws.Cells.Item(341, 7).Value = "Brócoli"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Segunda"
$ws.Cells.Item(341, 10).Value = 120
$ws.Cells.Item(341, 11).Value = 400
$ws.Cells.Item(341, 12).Value = 400
$ws.Cells.Item(341, 13).Value = 400
$ws.Cells.Item(341, 14).Value = "$/unidad"
$ws.Cells.Item(341, 15).Value = "Región del Maule"
$ws.Cells.Item(341, 16).Value = 400
$ws.Cells.Item(341, 17).Value = 1
$ws.Cells.Item(341, 18).Value = "Hortaliza"
